$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.040902506845782
$ws.Range("C2").Value = 0.4566974898799856
$ws.Range("D2").Value = 0.03342088562027357
$ws.Range("E2").Value = 0.0995373441219094
$ws.Range("F2").Value = 4.747147292256614
$ws.Range("J2").Value = 0.2477278073313016
$ws.Range("B3").Value = 1.94423878183062
$ws.Range("C3").Value = 0.4308694591482833
$ws.Range("D3").Value = 0.03323313197268973
$ws.Range("E3").Value = 0.09944323500833008
$ws.Range("F3").Value = 4.569403720894172
$ws.Range("J3").Value = 0.2409718732438932
$ws.Range("B4").Value = 1.886734367927659
$ws.Range("C4").Value = 0.4154811655082256
$ws.Range("D4").Value = 0.03314560468809802
$ws.Range("E4").Value = 0.09946038421876935
$ws.Range("F4").Value = 4.46210882643328
$ws.Range("J4").Value = 0.2369674581044166
$ws.Range("B5").Value = 1.863761877539275
$ws.Range("C5").Value = 0.4093274142138625
$ws.Range("D5").Value = 0.03311680360571501
$ws.Range("E5").Value = 0.09948619738362652
$ws.Range("F5").Value = 4.418837834462124
$ws.Range("J5").Value = 0.2353714024204976
$ws.Range("B6").Value = 1.859975056216285
$ws.Range("C6").Value = 0.4083126286184893
$ws.Range("D6").Value = 0.03311243270152531
$ws.Range("E6").Value = 0.09949162013105095
$ws.Range("F6").Value = 4.411679755150146
$ws.Range("J6").Value = 0.2351085286131109
$ws.Range("B7").Value = 1.886422690812594
$ws.Range("C7").Value = 0.4153977011966958
$ws.Range("D7").Value = 0.03314518861027338
$ws.Range("E7").Value = 0.09946065614486699
$ws.Range("F7").Value = 4.461523439782809
$ws.Range("J7").Value = 0.2369457887705124
$ws.Range("B8").Value = 2.007187204375441
$ws.Range("C8").Value = 0.4476936866680603
$ws.Range("D8").Value = 0.0333503199134455
$ws.Range("E8").Value = 0.09948933134741722
$ws.Range("F8").Value = 4.685473094683658
$ws.Range("J8").Value = 0.2453682961404979
$ws.Range("B9").Value = 2.258845136086506
$ws.Range("C9").Value = 0.5148137540723496
$ws.Range("D9").Value = 0.03397793590937681
$ws.Range("E9").Value = 0.1001412186588375
$ws.Range("F9").Value = 5.139728110102197
$ws.Range("J9").Value = 0.2630435676030629
$ws.Range("B10").Value = 2.453063006835407
$ws.Range("C10").Value = 0.566522460298529
$ws.Range("D10").Value = 0.03458366768109755
$ws.Range("E10").Value = 0.100985251394178
$ws.Range("F10").Value = 5.483382930874427
$ws.Range("J10").Value = 0.2767630222234203
$ws.Range("B11").Value = 2.543506417202082
$ws.Range("C11").Value = 0.5905860372107554
$ws.Range("D11").Value = 0.03489222983218809
$ws.Range("E11").Value = 0.1014489930893276
$ws.Range("F11").Value = 5.642031787625399
$ws.Range("J11").Value = 0.2831696934213142
$ws.Range("B12").Value = 2.578060842260982
$ws.Range("C12").Value = 0.5997776360853777
$ws.Range("D12").Value = 0.03501395315198863
$ws.Range("E12").Value = 0.1016361121986264
$ws.Range("F12").Value = 5.702453857129456
$ws.Range("J12").Value = 0.2856200257045316
$ws.Range("B13").Value = 2.570605268899897
$ws.Range("C13").Value = 0.5977945128639703
$ws.Range("D13").Value = 0.03498751882579398
$ws.Range("E13").Value = 0.1015953002135674
$ws.Range("F13").Value = 5.68942534101177
$ws.Range("J13").Value = 0.2850912163602146
$ws.Range("B14").Value = 2.54634308197825
$ws.Range("C14").Value = 0.5913406389537386
$ws.Range("D14").Value = 0.03490214557461258
$ws.Range("E14").Value = 0.1014641565882464
$ws.Range("F14").Value = 5.646995761261167
$ws.Range("J14").Value = 0.2833707948190494
$ws.Range("B15").Value = 2.531521710179788
$ws.Range("C15").Value = 0.5873978185240389
$ws.Range("D15").Value = 0.03485049110256
$ws.Range("E15").Value = 0.1013853274416725
$ws.Range("F15").Value = 5.621051725893665
$ws.Range("J15").Value = 0.282320160459264
$ws.Range("B16").Value = 2.44719481265372
$ws.Range("C16").Value = 0.5649608582941141
$ws.Range("D16").Value = 0.03456417780179066
$ws.Range("E16").Value = 0.1009565533943473
$ws.Range("F16").Value = 5.473062512912747
$ws.Range("J16").Value = 0.2763477020507139
$ws.Range("B17").Value = 2.396002170282657
$ws.Range("C17").Value = 0.5513360935780725
$ws.Range("D17").Value = 0.03439708108702888
$ws.Range("E17").Value = 0.1007139751532513
$ws.Range("F17").Value = 5.382878081510455
$ws.Range("J17").Value = 0.2727265280358182
$ws.Range("B18").Value = 2.366754213258275
$ws.Range("C18").Value = 0.5435503239976356
$ws.Range("D18").Value = 0.03430407075797604
$ws.Range("E18").Value = 0.1005819581147165
$ws.Range("F18").Value = 5.331223898052656
$ws.Range("J18").Value = 0.2706592862237329
$ws.Range("B19").Value = 2.356885004639878
$ws.Range("C19").Value = 0.5409228838314561
$ws.Range("D19").Value = 0.03427310742519296
$ws.Range("E19").Value = 0.1005385477622944
$ws.Range("F19").Value = 5.31377168197227
$ws.Range("J19").Value = 0.269962012908735
$ws.Range("B20").Value = 2.401431324698024
$ws.Range("C20").Value = 0.5527812015207019
$ws.Range("D20").Value = 0.03441454711558833
$ws.Range("E20").Value = 0.1007390207598569
$ws.Range("F20").Value = 5.39245577138891
$ws.Range("J20").Value = 0.2731103943978894
$ws.Range("B21").Value = 2.553461147161727
$ws.Range("C21").Value = 0.5932341357357132
$ws.Range("D21").Value = 0.03492708834140501
$ws.Range("E21").Value = 0.101502363958371
$ws.Range("F21").Value = 5.659448893999752
$ws.Range("J21").Value = 0.2838754621576669
$ws.Range("B22").Value = 2.654604143099732
$ws.Range("C22").Value = 0.6201350695388328
$ws.Range("D22").Value = 0.03529055585753582
$ws.Range("E22").Value = 0.1020683622279179
$ws.Range("F22").Value = 5.835961989752775
$ws.Range("J22").Value = 0.2910527067986663
$ws.Range("B23").Value = 2.600457565000795
$ws.Range("C23").Value = 0.6057347166834575
$ws.Range("D23").Value = 0.03509391564955422
$ws.Range("E23").Value = 0.1017601248553213
$ws.Range("F23").Value = 5.741565063850828
$ws.Range("J23").Value = 0.2872089641974469
$ws.Range("B24").Value = 2.398976233155622
$ws.Range("C24").Value = 0.5521277210335711
$ws.Range("D24").Value = 0.03440664121551151
$ws.Range("E24").Value = 0.1007276744562624
$ws.Range("F24").Value = 5.388125094268958
$ws.Range("J24").Value = 0.2729368029134918
$ws.Range("B25").Value = 2.189148089425487
$ws.Range("C25").Value = 0.4962417108757222
$ws.Range("D25").Value = 0.03378328206628112
$ws.Range("E25").Value = 0.09990090330498447
$ws.Range("F25").Value = 5.015151403091835
$ws.Range("J25").Value = 0.2581351263487903
